$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.300.42'
$ws.Range('E2').Value = '  -2.00%  '
$ws.Range('D3').Value = '3.605.64'
$ws.Range('E3').Value = '  -2.28%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''627.99'
$ws.Range('E5').Value = '  -6.44%  '
$ws.Range('D6').Value = '''156.70'
$ws.Range('E6').Value = '  -2.46%  '
$ws.Range('D7').Value = '3.604.56'
$ws.Range('E7').Value = '  -2.22%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '''0.490'
$ws.Range('E9').Value = '  -1.93%  '
$ws.Range('E10').Value = '  -2.73%  '
$ws.Range('D11').Value = '''7.01'
$ws.Range('E11').Value = '  -1.14%  '
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('E13').Value = '  -3.29%  '
$ws.Range('D14').Value = '4.213.27'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').Value = '''32.14'
$ws.Range('E15').Value = '  -3.20%  '
$ws.Range('D16').Value = '3.601.19'
$ws.Range('E16').Value = '  -2.35%  '
$ws.Range('D17').Value = '68.279.89'
$ws.Range('E17').Value = '  -1.98%  '
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').Value = '''6.45'
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('D20').Value = '''15.69'
$ws.Range('E20').Value = '  -2.73%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '''459.20'
$ws.Range('E21').Value = '  -2.63%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '''9.92'
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('D23').Value = '''0.644'
$ws.Range('E23').Value = '  -0.62%  '
$ws.Range('D24').Value = '''78.14'
$ws.Range('E24').Value = '  -2.17%  '
$ws.Range('D25').Value = '3.746.60'
$ws.Range('E25').Value = '  -2.38%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = '''10.81'
$ws.Range('E27').Value = '  -1.48%  '
$ws.Range('E28').Value = '  -7.84%  '
$ws.Range('D29').Value = '''8.48'
$ws.Range('E29').Value = '  -6.74%  '
$ws.Range('E30').Value = '  -3.27%  '
$ws.Range('D31').Value = '''1.65'
$ws.Range('E31').Value = '  -4.15%  '
$ws.Range('D32').Value = '''1.00'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('E33').Value = '  -4.45%  '
$ws.Range('D34').Value = '''26.17'
$ws.Range('E34').Value = '  -2.54%  '
$ws.Range('E35').Value = '  -4.58%  '
$ws.Range('D36').Value = '3.603.47'
$ws.Range('E36').Value = '  -2.29%  '
$ws.Range('D37').Value = '''6.25'
$ws.Range('E37').Value = '  -3.97%  '
$ws.Range('D38').Value = '''8.22'
$ws.Range('E38').Value = '  -3.42%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '''0.999'
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '''177.25'
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('D42').Value = '''5.67'
$ws.Range('E42').Value = '  -7.61%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''2.16'
$ws.Range('E43').Value = '  -4.61%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').Value = '''0.0886'
$ws.Range('E44').Value = '  -2.33%  '
$ws.Range('D45').Value = '''0.909'
$ws.Range('E45').Value = '  -2.88%  '
$ws.Range('D46').Value = '''29.09'
$ws.Range('E46').Value = '  +4.53%  '
$ws.Range('D47').Value = '''46.09'
$ws.Range('E48').Value = '  -5.09%  '
$ws.Range('D49').Value = '''7.75'
$ws.Range('E49').Value = '  -1.78%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = '''1.21'
$ws.Range('E50').Value = '  -6.58%  '
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').Value = '''1.02'
$ws.Range('E51').Value = '  -5.88%  '
